$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 596-613: columns F (6) and G (7)
$updates = @(
    @{Row=596; F=29201; G=943},
    @{Row=597; F=29545; G=954},
    @{Row=598; F=15474; G=706},
    @{Row=599; F=16598; G=872},
    @{Row=600; F=39977; G=1675},
    @{Row=601; F=31671; G=1332},
    @{Row=602; F=30068; G=1290},
    @{Row=603; F=31881; G=1523},
    @{Row=604; F=29903; G=1527},
    @{Row=605; F=14799; G=1024},
    @{Row=606; F=14274; G=1262},
    @{Row=607; F=10821; G=953},
    @{Row=608; F=45767; G=2890},
    @{Row=609; F=36318; G=2153},
    @{Row=610; F=33474; G=1906},
    @{Row=611; F=33830; G=2107},
    @{Row=612; F=15831; G=1400},
    @{Row=613; F=21662; G=1893}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 6).Value = $u.F
    $ws.Cells.Item($u.Row, 7).Value = $u.G
}

# Add new rows 614-616
$newRows = @(
    @{Row=614; A=44508; B=526608; C=19537; D=4958; E=13314; F=46809; G=3284},
    @{Row=615; A=44509; B=533663; C=24791; D=7055; E=13367; F=35799; G=2291},
    @{Row=616; A=44510; B=540209; C=22967; D=6546; E=13405; F=27186; G=1859}
)

foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 1).NumberFormat = "yyyy-mm-dd"
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 6).Value = $r.F
    $ws.Cells.Item($r.Row, 7).Value = $r.G
}
